$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.299
$ws.Range("C2").Value = 6.6536
$ws.Range("D2").Value = 155.21038174364
$ws.Range("E2").Value = 3.930683045250641

$ws.Range("B3").Value = 9.481299999999999
$ws.Range("C3").Value = 11.435
$ws.Range("D3").Value = 256.47615631891
$ws.Range("E3").Value = 677430000
